$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the precise timestamp value on row 15 (floating point precision correction)
$ws.Range("A15").Value = 45866.66691041666

# Add new row 16 data
$ws.Range("A16").Value = 45866.7086142748
$ws.Range("B16").Value = 2025
$ws.Range("C16").Value = 31
$ws.Range("D16").Value = 20.11
$ws.Range("E16").Value = 72.67
$ws.Range("F16").Value = 119.53
$ws.Range("G16").Value = 9.35
$ws.Range("H16").Value = "ESE"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = "17:00:24"

# Match number format of column A (date/time formatting) from row 15
$ws.Range("A16").NumberFormat = $ws.Range("A15").NumberFormat
